# AddressTypes.xlsx catalogue reload
#  - rename the data sheet
#  - add PriorityForOKATO / Category columns
#  - fix a typo in the "By passport" label
#  - add priority/category values to the existing rows
#  - add a new "By talon" address type row
#
# Numeric-looking / boolean-looking text ("1", "2", "True", "False", ...) has
# to stay stored as plain text (matches the original workbook, which keeps
# everything as shared strings), but a direct `.Value = "1"` assignment gets
# auto-coerced by Excel into a real number/boolean. Routing the text through
# a helper cell's `="..."` formula and pasting it back as values keeps the
# literal text *and* avoids bolting a new number-format style onto the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        $Sheet,
        [string]$CellRef,
        [string]$Text
    )
    $helper = $Sheet.Range("ZZ1")
    $helper.Formula = '=""' + $Text + '""'
    $helper.Copy()
    $Sheet.Range($CellRef).PasteSpecial(-4163)
    $helper.ClearContents()
}

# 1. Rename the worksheet
$ws.Name = "Data"

# 2. New header cells (plain text, safe to assign directly)
$ws.Range("D1").Value = "PriorityForOKATO"
$ws.Range("E1").Value = "Category"

# 3. Row 2 (Id = 1): fix the typo, add priority + category
$ws.Range("B2").Value = "По паспарту"
Set-TextValue $ws "D2" "2"
$ws.Range("E2").Value = "|Registry|"

# 4. Row 3 (Id = 2): add priority + category
Set-TextValue $ws "D3" "1"
$ws.Range("E3").Value = "|Registry|"

# 5. Row 4 (Id = 3): add priority + category
Set-TextValue $ws "D4" "3"
$ws.Range("E4").Value = "|Registry|"

# 6. Row 5 (new, Id = 4): "By talon" address type
Set-TextValue $ws "A5" "4"
$ws.Range("B5").Value = "По талону"
Set-TextValue $ws "C5" "True"
Set-TextValue $ws "D5" "4"
$ws.Range("E5").Value = "|Talon|"
